$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.308.39'
$ws.Range('E2').Value = '  -1.00%  '

$ws.Range('D3').Value = '1.860.71'
$ws.Range('E3').Value = '  -1.17%  '

$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '233.79'
$ws.Range('E5').Value = '  -2.34%  '

$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  +0.05%  '

$ws.Range('D7').Value = '0.4763'
$ws.Range('E7').Value = '  -0.97%  '

$ws.Range('D8').Value = '0.2754'
$ws.Range('E8').Value = '  -2.85%  '

$ws.Range('D9').Value = '0.06446'
$ws.Range('E9').Value = '  -1.60%  '

$ws.Range('D10').Value = '1.912.01'
$ws.Range('E10').Value = '  +1.48%  '

$ws.Range('D11').Value = '0.07435'
$ws.Range('E11').Value = '  -0.73%  '

$ws.Range('D12').Value = '16.13'
$ws.Range('E12').Value = '  -3.52%  '

$ws.Range('D13').Value = '5.000'
$ws.Range('E13').Value = '  -2.07%  '

$ws.Range('D14').Value = '86.00'
$ws.Range('E14').Value = '  -2.98%  '

$ws.Range('D15').Value = '0.6330'
$ws.Range('E15').Value = '  -5.05%  '

$ws.Range('D16').Value = '30.279.14'
$ws.Range('E16').Value = '  -0.89%  '

$ws.Range('D17').Value = '0.9999'
$ws.Range('E17').Value = '  +0.11%  '

$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = '231.96'
$ws.Range('E18').Value = '  +3.32%  '

$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '12.83'
$ws.Range('E19').Value = '  -4.12%  '

$ws.Range('D20').Value = '0.000007381'
$ws.Range('E20').Value = '  -3.26%  '

$ws.Range('D21').Value = '2.097.29'
$ws.Range('E21').Value = '  -1.99%  '

$ws.Range('D22').Value = '0.9998'
$ws.Range('E22').Value = '  +0.13%  '

$ws.Range('D23').Value = '5.108'
$ws.Range('E23').Value = '  -4.48%  '

$ws.Range('B24').Value = 'BitDAO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D24').Value = '0.3935'
$ws.Range('E24').Value = '  -1.31%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = '6.016'
$ws.Range('E25').Value = '  -3.40%  '

$ws.Range('D26').Value = '9.297'
$ws.Range('E26').Value = '  -0.82%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '167.73'
$ws.Range('E27').Value = '  +0.58%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '17.92'
$ws.Range('E28').Value = '  -3.85%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '1.860'
$ws.Range('E29').Value = '  -6.03%  '

$ws.Range('E30').Value = '  -4.91%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.1008'
$ws.Range('E31').Value = '  +6.46%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.232'
$ws.Range('E32').Value = '  -2.50%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '3.921'
$ws.Range('E33').Value = '  -2.96%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.04896'
$ws.Range('E34').Value = '  -2.86%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.149'
$ws.Range('E35').Value = '  -5.03%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7247'
$ws.Range('E36').Value = '  -3.39%  '

$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = '0.9992'
$ws.Range('E37').Value = '  +0.19%  '

$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.691'
$ws.Range('E38').Value = '  -0.26%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01952'
$ws.Range('E39').Value = '  +6.14%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.630'
$ws.Range('E40').Value = '  +0.33%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '0.9098'
$ws.Range('E41').Value = '  +0.13%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '1.989'
$ws.Range('E42').Value = '  -5.06%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '105.71'
$ws.Range('E43').Value = '  -0.47%  '

$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '0.9997'
$ws.Range('E44').Value = '  -0.66%  '

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.4124'
$ws.Range('E45').Value = '  -4.09%  '

$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '5.557'
$ws.Range('E46').Value = '  -4.87%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '7.072'
$ws.Range('E47').Value = '  -5.82%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '61.33'
$ws.Range('E48').Value = '  -5.47%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.1209'
$ws.Range('E49').Value = '  -5.96%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.780'
$ws.Range('E50').Value = '  -1.75%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.404'
$ws.Range('E51').Value = '  -5.21%  '
